$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.352.42'
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('D3').Value = '2.518.06'
$ws.Range('E3').Value = '  +3.31%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '542.13'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('D6').Value = '144.76'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = '2.550.67'
$ws.Range('E9').Value = '  +3.96%  '
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '5.60'
$ws.Range('E12').Value = '  +5.75%  '
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = '2.963.86'
$ws.Range('E14').Value = '  +2.99%  '
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '59.279.30'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '2.544.61'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').Value = '11.26'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '326.61'
$ws.Range('E21').Value = '  +1.10%  '
$ws.Range('E22').Value = '  +3.36%  '
$ws.Range('D23').Value = '5.85'
$ws.Range('E23').Value = '  +3.22%  '
$ws.Range('D24').Value = '62.18'
$ws.Range('E24').Value = '  +2.69%  '
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('D26').Value = '0.165'
$ws.Range('E26').Value = '  +3.23%  '
$ws.Range('D27').Value = '0.992'
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('D28').Value = '8.01'
$ws.Range('E28').Value = '  +4.97%  '
$ws.Range('E29').Value = '  +3.78%  '
$ws.Range('D30').Value = '0.0₃0786'
$ws.Range('E30').Value = '  +2.76%  '
$ws.Range('D31').Value = '1.83'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('E33').Value = '  +10.52%  '
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').Value = '157.14'
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').Value = '18.71'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').Value = '4.40'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('E38').Value = '  -3.25%  '
$ws.Range('D39').Value = '5.68'
$ws.Range('E39').Value = '  -1.96%  '
$ws.Range('D41').Value = '299.74'
$ws.Range('E41').Value = '  -3.62%  '
$ws.Range('D42').Value = '3.71'
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').Value = '0.831'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('E45').Value = '  +4.66%  '
$ws.Range('D46').Value = '10.80'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('E48').Value = '  +2.64%  '
$ws.Range('D49').Value = '122.98'
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('B50').Value = 'Hedera'
$ws.Range('C50').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D50').Value = '0.0517'
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '0.0229'
$ws.Range('E51').Value = '  +0.02%  '
